$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: wipe existing cell values (keep styles attached to A1/B1/A6/B7:B13) ---
$ws.Cells.ClearContents()

# --- Step 2: restructure the <cols> metadata so the two custom-width column
# spans grow from (1-4 / 2-13) to (1-10 / 2-17), matching the widened table
# (columns for teams B and C, plus the Q columns). Do this while the data
# area is empty so the insert/delete doesn't scramble real content. ---
$ws.Columns("C:H").Insert()
$ws.Columns("L:M").Delete()

# Remove the leftover blank/format noise the column insert left behind in
# C1:H13 (it inherited formatting from neighbouring cells).
$ws.Range("C1:H13").Clear()

# --- Step 3: push the weekly-schedule block (row 6 header + the 7 weekday
# rows 7-13) down by 4 rows, to rows 10 and 11-17, making room for the
# expanded roster block above it. This carries A6's and B7:B13's styles
# along automatically. ---
$ws.Rows("6:9").Insert()

# --- Step 4: fill in the roster block (rows 1-7) ---
$ws.Range("A1").Value = "Equipes"
$ws.Range("B1").Value = "A"

$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "C"

$ws.Range("B2").Value = "F. Mayweather"
$ws.Range("B3").Value = "C. Ronaldo"
$ws.Range("B4").Value = "L. Messi"
$ws.Range("B5").Value = "K. Bryant"
$ws.Range("B6").Value = "R. Federer"
$ws.Range("B7").Value = "P. Mickelson"

$ws.Range("C2").Value = "R. Nadal"
$ws.Range("C3").Value = "M. Ryan"
$ws.Range("C4").Value = "M. Pacquiao"
$ws.Range("C5").Value = "Z. Ibrahimović"
$ws.Range("C6").Value = "D. Rose"
$ws.Range("C7").Value = "G. Bale"

$ws.Range("D2").Value = "R. Falcao"
$ws.Range("D3").Value = "M. Özil"
$ws.Range("D4").Value = "N. Djokovic"
$ws.Range("D5").Value = "M. Stafford"
$ws.Range("D6").Value = "L. Hamilton"
$ws.Range("D7").Value = "K. Durant"

# --- Step 5: fill in the weekly-schedule header (row 10) ---
$ws.Range("A10").Value = "Semaine 13"
$ws.Range("B10").Value = "Horaire"
$ws.Range("C10").Value = "Q1"
$ws.Range("D10").Value = "Q2"
$ws.Range("E10").Value = "Q3"

# --- Step 6: fill in the weekday rows (11-17), already carrying style 2 on
# column B from the row-insert shift; fill in the Q1/Q2 matchups ---
$ws.Range("B11").Value = "Lundi`n2022-03-28"
$ws.Range("B12").Value = "Mardi`n2022-03-29"
$ws.Range("B13").Value = "Mercredi`n2022-03-30"
$ws.Range("B14").Value = "Jeudi`n2022-03-31"
$ws.Range("B15").Value = "Vendredi`n2022-04-01"
$ws.Range("B16").Value = "Samedi`n2022-04-02"
$ws.Range("B17").Value = "dimanche`n2022-04-03"

$ws.Range("C11").Value = "A B"
$ws.Range("D11").Value = "C A"
$ws.Range("C12").Value = "B C"
$ws.Range("D12").Value = "A B"
$ws.Range("C13").Value = "C A"
$ws.Range("D13").Value = "B C"
$ws.Range("C14").Value = "A B"
$ws.Range("D14").Value = "C A"
$ws.Range("C15").Value = "B C"
$ws.Range("D15").Value = "A"

# Multi-line day values trigger an auto row-height calc on assignment;
# AutoFit() clears the resulting ht/customHeight markup back to default.
$ws.Rows("11:17").AutoFit()

# --- Step 7: footer note (row 20), styled like the other section titles ---
$ws.Range("A1").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = "Modele : repartition concentre h-pers = 900"
